$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 27; this shifts the existing rows 27-54 down to 28-55.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record.
$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(27, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(27, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(27, 4).Value = 44614
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(27, 6).Value = 100112009
$ws.Cells.Item(27, 7).Value = "Acelga"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 300
$ws.Cells.Item(27, 11).Value = 1200
$ws.Cells.Item(27, 12).Value = 1500
$ws.Cells.Item(27, 13).Value = 1350
$ws.Cells.Item(27, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 450
$ws.Cells.Item(27, 17).Value = 3
$ws.Cells.Item(27, 18).Value = "Hortaliza"
